# MTP049 COP rows delivered
#
# The "New Priority" header note that was parked in D1 (row 1) moves down
# into its own bold + underlined section label in A12, ahead of where the
# MTP049-onwards rows will be delivered, and its text is expanded.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# D1 previously held the "New Priority" note; clear its content (style stays).
$ws.Range("D1").Value = ""

# New section-header cell, two rows below the first table, carrying the
# (expanded) text that used to live in D1.
$ws.Range("A12").Value = "New Priorities MTP049 onwards"
$ws.Range("A12").Font.Bold = $true
$ws.Range("A12").Font.Underline = $true
$ws.Range("A12").WrapText = $true

# Row 1 no longer needs the extra height that wrapping "New Priority" in D1
# required - let it size back to the default.
$ws.Rows.Item(1).AutoFit()

# Reflect the new point of interest as the active selection.
[void]$ws.Range("A12").Select()

Write-Host "done"
